$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price report row was inserted above the existing row 24,
# pushing the previous rows 24-32 down to 25-33 (dimension grows to R33).
$ws.Rows.Item(24).EntireRow.Insert()

# Populate the newly inserted row 24 with this week's data.
$ws.Cells.Item(24, 1).Value = 8
$ws.Cells.Item(24, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(24, 3).Value = "Coquimbo"
$ws.Cells.Item(24, 4).Value = 44846
$ws.Cells.Item(24, 5).Value = 4
$ws.Cells.Item(24, 6).Value = 100112026
$ws.Cells.Item(24, 7).Value = "Haba"
$ws.Cells.Item(24, 8).Value = "Sin especificar"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 600
$ws.Cells.Item(24, 11).Value = 7500
$ws.Cells.Item(24, 12).Value = 8000
$ws.Cells.Item(24, 13).Value = 7750
$ws.Cells.Item(24, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(24, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(24, 16).Value = 310
$ws.Cells.Item(24, 17).Value = 25
$ws.Cells.Item(24, 18).Value = "Hortaliza"
